$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: "Value" (uoa totals) becomes "Year"; loses its old centered style ---
$ws.Range("B1").Value = "Year"
$ws.Range("B1").ClearFormats()

$ws.Range("B2").Value = 2007
$ws.Range("B2").ClearFormats()

$ws.Range("B3").Value = 2014
$ws.Range("B3").ClearFormats()

# old row4 ("uoa3" / 20000) is no longer part of the table -> wipe it out completely
$ws.Range("A4").Clear()
$ws.Range("B4").Clear()

# --- New column C: Currency ---
$ws.Range("C1").Value = "Currency"
$ws.Range("C2").Value = "USD"
$ws.Range("C3").Value = "USD"

# --- New column D: Value (moved from column B, same look & feel as the old column) ---
$ws.Range("D1").Value = "Value"
$ws.Range("D1").HorizontalAlignment = -4108

$ws.Range("D2").Value = 16823445.68
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").NumberFormat = "#,##0.00"

$ws.Range("D3").Value = 4005582.31
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").NumberFormat = "#,##0.00"

# keep the numeric look going for the (still empty) row below the table
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").NumberFormat = "#,##0.00"

$ws.Columns.Item(4).ColumnWidth = 11.9

$ws.Range("A3").Select()
